$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1921
$ws1.Range("F5").Value = 864
$ws1.Range("F6").Value = 269

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1921
$ws4.Range("F6").Value = 864
$ws4.Range("F7").Value = 269
